$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2119.0588
$ws.Range("I15").Value = 2119.0588
$ws.Range("K15").Value = 6357.176399999999
$ws.Range("M15").Value = -6188.176399999999
$ws.Range("H69").Value = 22333.334
$ws.Range("I69").Value = 14000
$ws.Range("J69").Value = 23090.908
$ws.Range("K69").Value = 42000
$ws.Range("L69").Value = 69272.724
$ws.Range("M69").Value = -41126
$ws.Range("N69").Value = -71020.724
$ws.Range("H72").Value = 22333.334
$ws.Range("I72").Value = 14000
$ws.Range("J72").Value = 23090.908
$ws.Range("K72").Value = 126000
$ws.Range("L72").Value = 207818.172
$ws.Range("M72").Value = -121632
$ws.Range("N72").Value = -216554.172
$ws.Range("H92").Value = 629.05884
$ws.Range("I92").Value = 532.7857
$ws.Range("K92").Value = 532.7857
$ws.Range("M92").Value = 715.2143
$ws.Range("H127").Value = 2044.8334
$ws.Range("I127").Value = 1826.6666
$ws.Range("J127").Value = 2699.3333
$ws.Range("K127").Value = 5479.9998
$ws.Range("L127").Value = 8097.999899999999
$ws.Range("M127").Value = -519.9997999999996
$ws.Range("N127").Value = -18017.9999
$ws.Range("H129").Value = 1719
$ws.Range("J129").Value = 2698.5
$ws.Range("L129").Value = 8095.5
$ws.Range("N129").Value = -18095.5
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H138").Value = 5645.4736
$ws.Range("J138").Value = 5710.4375
$ws.Range("L138").Value = 17131.3125
$ws.Range("N138").Value = -27411.3125

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 1151.5
$ws.Range("I28").Value = 1151.5
$ws.Range("K28").Value = 1151.5
$ws.Range("M28").Value = -959.5
$ws.Range("H99").Value = 1151.5
$ws.Range("I99").Value = 1151.5
$ws.Range("K99").Value = 1151.5
$ws.Range("M99").Value = 1843.5
$ws.Range("H122").Value = 13404.223
$ws.Range("I122").Value = 21680.2
$ws.Range("J122").Value = 3059.25
$ws.Range("K122").Value = 65040.60000000001
$ws.Range("L122").Value = 9177.75
$ws.Range("M122").Value = -62590.60000000001
$ws.Range("N122").Value = -14077.75
$ws.Range("H132").Value = 1099.5
$ws.Range("I132").Value = 1074.25
$ws.Range("K132").Value = 3222.75
$ws.Range("M132").Value = -692.75

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1444.5
$ws.Range("I86").Value = 1444.5
$ws.Range("K86").Value = 1444.5
$ws.Range("M86").Value = -321.5
$ws.Range("H89").Value = 1444.5
$ws.Range("I89").Value = 1444.5
$ws.Range("K89").Value = 7222.5
$ws.Range("M89").Value = -1606.5
$ws.Range("H94").Value = 1440.9286
$ws.Range("I94").Value = 1406.75
$ws.Range("K94").Value = 1406.75
$ws.Range("M94").Value = -955.75

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2150.1667
$ws.Range("I31").Value = 1600.0834
$ws.Range("J31").Value = 2516.889
$ws.Range("K31").Value = 1600.0834
$ws.Range("L31").Value = 2516.889
$ws.Range("M31").Value = -1305.0834
$ws.Range("N31").Value = -3106.889
$ws.Range("H34").Value = 2150.1667
$ws.Range("I34").Value = 1600.0834
$ws.Range("J34").Value = 2516.889
$ws.Range("K34").Value = 1600.0834
$ws.Range("L34").Value = 2516.889
$ws.Range("M34").Value = -1398.0834
$ws.Range("N34").Value = -2920.889
$ws.Range("H86").Value = 3266.6667
$ws.Range("I86").Value = 3260.2
$ws.Range("K86").Value = 3260.2
$ws.Range("M86").Value = -2137.2
$ws.Range("H89").Value = 3266.6667
$ws.Range("I89").Value = 3260.2
$ws.Range("K89").Value = 16301
$ws.Range("M89").Value = -10685
$ws.Range("H94").Value = 4037.0908
$ws.Range("I94").Value = 3518
$ws.Range("J94").Value = 4469.6665
$ws.Range("K94").Value = 3518
$ws.Range("L94").Value = 4469.6665
$ws.Range("M94").Value = -3067
$ws.Range("N94").Value = -5371.6665
$ws.Range("H124").Value = 89872.5
$ws.Range("J124").Value = 89872.5
$ws.Range("L124").Value = 89872.5
$ws.Range("N124").Value = -94782.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 3608.6
$ws.Range("J68").Value = 3608.6
$ws.Range("L68").Value = 10825.8
$ws.Range("N68").Value = -12447.8
$ws.Range("H71").Value = 3608.6
$ws.Range("J71").Value = 3608.6
$ws.Range("L71").Value = 32477.4
$ws.Range("N71").Value = -40589.39999999999
$ws.Range("H114").Value = 742.3333
$ws.Range("I114").Value = 742.3333
$ws.Range("K114").Value = 2226.9999
$ws.Range("M114").Value = 1027.0001
$ws.Range("H121").Value = 3647.7058
$ws.Range("I121").Value = 1043.5
$ws.Range("K121").Value = 3130.5
$ws.Range("M121").Value = -1820.5
$ws.Range("H131").Value = 2095.1052
$ws.Range("I131").Value = 1016.4286
$ws.Range("J131").Value = 2724.3333
$ws.Range("K131").Value = 3049.2858
$ws.Range("L131").Value = 8172.999899999999
$ws.Range("M131").Value = 1990.7142
$ws.Range("N131").Value = -18252.9999
$ws.Range("H134").Value = 9670.385
$ws.Range("I134").Value = 14139.625
$ws.Range("K134").Value = 42418.875
$ws.Range("M134").Value = -37348.875
$ws.Range("H137").Value = 2562
$ws.Range("I137").Value = 2042.25
$ws.Range("J137").Value = 2977.8
$ws.Range("K137").Value = 6126.75
$ws.Range("L137").Value = 8933.400000000001
$ws.Range("M137").Value = -1026.75
$ws.Range("N137").Value = -19133.4
$ws.Range("H140").Value = 60878.21
$ws.Range("I140").Value = 111909.78
$ws.Range("J140").Value = 14949.8
$ws.Range("K140").Value = 335729.34
$ws.Range("L140").Value = 44849.39999999999
$ws.Range("M140").Value = -330549.34
$ws.Range("N140").Value = -55209.39999999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2779.0527
$ws.Range("I122").Value = 2856.0557
$ws.Range("J122").Value = 1393
$ws.Range("K122").Value = 8568.167099999999
$ws.Range("L122").Value = 4179
$ws.Range("M122").Value = -6118.167099999999
$ws.Range("N122").Value = -9079

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1933.5454
$ws.Range("I68").Value = 1561.5
$ws.Range("K68").Value = 1561.5
$ws.Range("M68").Value = -812.5
$ws.Range("H71").Value = 1933.5454
$ws.Range("I71").Value = 1561.5
$ws.Range("K71").Value = 7807.5
$ws.Range("M71").Value = -4063.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1562.6
$ws.Range("I107").Value = 1549.6
$ws.Range("K107").Value = 4648.799999999999
$ws.Range("M107").Value = -2728.799999999999
$ws.Range("H132").Value = 1172.25
$ws.Range("I132").Value = 1029.6666
$ws.Range("K132").Value = 3088.9998
$ws.Range("M132").Value = -558.9998000000001
